$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 15.32737633333333
$ws.Range("H2").Value = 45.982129
$ws.Range("I2").Value = 0.336907232375371
$ws.Range("J2").Value = 0.336907232375371
$ws.Range("M2").Value = 15.32737633333333
$ws.Range("N2").Value = 45.982129
$ws.Range("O2").Value = 0.336907232375371
$ws.Range("P2").Value = 0.336907232375371
$ws.Range("Q2").Value = 234.9284652636268
$ws.Range("R2").Value = 2114.356187372641
$ws.Range("S2").Value = 0.1135064832268322
$ws.Range("T2").Value = 0.1135064832268322
$ws.Range("G3").Value = 15.32737633333333
$ws.Range("H3").Value = 45.982129
$ws.Range("I3").Value = 0.336907232375371
$ws.Range("J3").Value = 0.336907232375371
$ws.Range("N3").Value = 7.967559000000001
$ws.Range("O3").Value = 0.05837764170244223
$ws.Range("P3").Value = 0.05837764170244223
$ws.Range("Q3").Value = 40.70725841701233
$ws.Range("R3").Value = 366.3653257531111
$ws.Range("S3").Value = 0.01966784969857085
$ws.Range("T3").Value = 0.01966784969857085
$ws.Range("G4").Value = 15.32737633333333
$ws.Range("H4").Value = 45.982129
$ws.Range("I4").Value = 0.336907232375371
$ws.Range("J4").Value = 0.336907232375371
$ws.Range("O4").Value = 0.6047151259221868
$ws.Range("P4").Value = 0.6047151259221868
$ws.Range("Q4").Value = 421.6733355736223
$ws.Range("R4").Value = 3795.060020162601
$ws.Range("S4").Value = 0.2037328994499679
$ws.Range("T4").Value = 0.2037328994499679
$ws.Range("H5").Value = 7.967559000000001
$ws.Range("I5").Value = 0.05837764170244223
$ws.Range("J5").Value = 0.05837764170244223
$ws.Range("M5").Value = 15.32737633333333
$ws.Range("N5").Value = 45.982129
$ws.Range("O5").Value = 0.336907232375371
$ws.Range("P5").Value = 0.336907232375371
$ws.Range("Q5").Value = 40.70725841701233
$ws.Range("R5").Value = 366.3653257531111
$ws.Range("S5").Value = 0.01966784969857085
$ws.Range("T5").Value = 0.01966784969857085
$ws.Range("H6").Value = 7.967559000000001
$ws.Range("I6").Value = 0.05837764170244223
$ws.Range("J6").Value = 0.05837764170244223
$ws.Range("N6").Value = 7.967559000000001
$ws.Range("O6").Value = 0.05837764170244223
$ws.Range("P6").Value = 0.05837764170244223
$ws.Range("R6").Value = 63.48199641848101
$ws.Range("S6").Value = 0.003407949050738722
$ws.Range("T6").Value = 0.003407949050738722
$ws.Range("H7").Value = 7.967559000000001
$ws.Range("I7").Value = 0.05837764170244223
$ws.Range("J7").Value = 0.05837764170244223
$ws.Range("O7").Value = 0.6047151259221868
$ws.Range("P7").Value = 0.6047151259221868
$ws.Range("R7").Value = 657.589486976271
$ws.Range("S7").Value = 0.03530184295313266
$ws.Range("T7").Value = 0.03530184295313266
$ws.Range("I8").Value = 0.6047151259221868
$ws.Range("J8").Value = 0.6047151259221868
$ws.Range("M8").Value = 15.32737633333333
$ws.Range("N8").Value = 45.982129
$ws.Range("O8").Value = 0.336907232375371
$ws.Range("P8").Value = 0.336907232375371
$ws.Range("Q8").Value = 421.6733355736223
$ws.Range("R8").Value = 3795.060020162601
$ws.Range("S8").Value = 0.2037328994499679
$ws.Range("T8").Value = 0.2037328994499679
$ws.Range("I9").Value = 0.6047151259221868
$ws.Range("J9").Value = 0.6047151259221868
$ws.Range("N9").Value = 7.967559000000001
$ws.Range("O9").Value = 0.05837764170244223
$ws.Range("P9").Value = 0.05837764170244223
$ws.Range("R9").Value = 657.589486976271
$ws.Range("S9").Value = 0.03530184295313266
$ws.Range("T9").Value = 0.03530184295313266
$ws.Range("I10").Value = 0.6047151259221868
$ws.Range("J10").Value = 0.6047151259221868
$ws.Range("O10").Value = 0.6047151259221868
$ws.Range("P10").Value = 0.6047151259221868
$ws.Range("S10").Value = 0.3656803835190862
$ws.Range("T10").Value = 0.3656803835190862
